$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "plan de vuelo en estructuras": replace the textual Airport/IAP/Runway/Waypoint
# lookups (KSEA / RNAVZRWY34L / 34L / BAKMN-JALON-RW34L) with numeric structure
# ids referencing the new "estructuras" (structures) table.
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(2, 3).Value = 14
$ws.Cells.Item(2, 4).Value = 36

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 14
$ws.Cells.Item(3, 4).Value = 37

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 5
$ws.Cells.Item(4, 3).Value = 14
$ws.Cells.Item(4, 4).Value = 38

# Update the active selection left by the editor
$ws.Range("D6").Select()
